$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 58
$ws.Range("H58").Value = 2935
$ws.Range("J58").Value = 6318.75
$ws.Range("L58").Value = 18956.25
$ws.Range("N58").Value = -19256.25
# Row 70
$ws.Range("H70").Value = 5375
$ws.Range("J70").Value = 7000
$ws.Range("L70").Value = 21000
$ws.Range("N70").Value = -21540
# Row 73
$ws.Range("H73").Value = 5375
$ws.Range("J73").Value = 7000
$ws.Range("L73").Value = 21000
$ws.Range("N73").Value = -22872
# Row 80
$ws.Range("H80").Value = 1972.3572
$ws.Range("I80").Value = 1941.2858
$ws.Range("K80").Value = 5823.857400000001
$ws.Range("M80").Value = -4825.857400000001
# Row 83
$ws.Range("H83").Value = 1972.3572
$ws.Range("I83").Value = 1941.2858
$ws.Range("K83").Value = 17471.5722
$ws.Range("M83").Value = -12479.5722
# Row 96
$ws.Range("H96").Value = 2209
$ws.Range("I96").Value = 2806.5
$ws.Range("J96").Value = 1014
$ws.Range("K96").Value = 8419.5
$ws.Range("L96").Value = 3042
$ws.Range("M96").Value = -7046.5
$ws.Range("N96").Value = -5788
# Row 100
$ws.Range("H100").Value = 2553.2727
$ws.Range("I100").Value = 2553.2727
$ws.Range("K100").Value = 2553.2727
$ws.Range("M100").Value = -2012.2727
# Row 132
$ws.Range("H132").Value = 76929400
$ws.Range("I132").Value = 76929400
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 230788200
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -230785670

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1939.7
$ws.Range("I61").Value = 1710.7778
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1710.7778
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -1498.7778
$ws.Range("N61").Value = -4424
# Row 97
$ws.Range("H97").Value = 2345.25
$ws.Range("I97").Value = 1287.4286
$ws.Range("J97").Value = 9750
$ws.Range("K97").Value = 1287.4286
$ws.Range("L97").Value = 9750
$ws.Range("M97").Value = -791.4286
$ws.Range("N97").Value = -10742
# Row 110
$ws.Range("H110").Value = 33642
$ws.Range("I110").Value = 22277.5
$ws.Range("K110").Value = 22277.5
$ws.Range("M110").Value = -20232.5
# Row 136
$ws.Range("H136").Value = 1939.7
$ws.Range("I136").Value = 1710.7778
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 5132.3334
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -2582.3334
$ws.Range("N136").Value = -17100
# Row 138
$ws.Range("H138").Value = 44999.875
$ws.Range("J138").Value = 44999.875
$ws.Range("L138").Value = 44999.875
$ws.Range("N138").Value = -55279.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1333.1
$ws.Range("I99").Value = 853.875
$ws.Range("K99").Value = 853.875
$ws.Range("M99").Value = 644.125
# Row 134
$ws.Range("H134").Value = 1679.4286
$ws.Range("I134").Value = 1679.4286
$ws.Range("K134").Value = 5038.2858
$ws.Range("M134").Value = -2503.2858

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 96
$ws.Range("H96").Value = 15896.111
$ws.Range("J96").Value = 15896.111
$ws.Range("L96").Value = 15896.111
$ws.Range("N96").Value = -21388.111
# Row 99
$ws.Range("H99").Value = 1639.1538
$ws.Range("I99").Value = 1413.5714
$ws.Range("J99").Value = 1902.3334
$ws.Range("K99").Value = 1413.5714
$ws.Range("L99").Value = 1902.3334
$ws.Range("M99").Value = 84.42859999999996
$ws.Range("N99").Value = -4898.3334
# Row 126
$ws.Range("H126").Value = 1639.1538
$ws.Range("I126").Value = 1413.5714
$ws.Range("J126").Value = 1902.3334
$ws.Range("K126").Value = 4240.7142
$ws.Range("L126").Value = 5707.0002
$ws.Range("M126").Value = -1770.7142
$ws.Range("N126").Value = -10647.0002
# Row 132
$ws.Range("H132").Value = 7778.1113
$ws.Range("I132").Value = 6875.375
$ws.Range("K132").Value = 20626.125
$ws.Range("M132").Value = -18096.125
# Row 134
$ws.Range("H134").Value = 1926
$ws.Range("I134").Value = 1926
$ws.Range("K134").Value = 5778
$ws.Range("M134").Value = -3243

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 8489.666999999999
$ws.Range("I80").Value = 8750
$ws.Range("J80").Value = 8229.333000000001
$ws.Range("K80").Value = 8750
$ws.Range("L80").Value = 8229.333000000001
$ws.Range("M80").Value = -7752
$ws.Range("N80").Value = -10225.333
# Row 83
$ws.Range("H83").Value = 8489.666999999999
$ws.Range("I83").Value = 8750
$ws.Range("J83").Value = 8229.333000000001
$ws.Range("K83").Value = 43750
$ws.Range("L83").Value = 41146.665
$ws.Range("M83").Value = -38758
$ws.Range("N83").Value = -51130.665
# Row 95
$ws.Range("H95").Value = 23669.5
$ws.Range("J95").Value = 23669.5
$ws.Range("L95").Value = 23669.5
$ws.Range("N95").Value = -29161.5
# Row 97
$ws.Range("H97").Value = 797.75
$ws.Range("I97").Value = 233.90909
$ws.Range("K97").Value = 233.90909
$ws.Range("M97").Value = 262.09091
# Row 107
$ws.Range("H107").Value = 1279.8
$ws.Range("I107").Value = 311
$ws.Range("J107").Value = 9999
$ws.Range("K107").Value = 311
$ws.Range("L107").Value = 9999
$ws.Range("M107").Value = 1609
$ws.Range("N107").Value = -13839
# Row 113
$ws.Range("H113").Value = 5800
$ws.Range("J113").Value = 7500
$ws.Range("L113").Value = 7500
$ws.Range("N113").Value = -11840
# Row 132
$ws.Range("H132").Value = 1722.7142
$ws.Range("I132").Value = 1447.5385
$ws.Range("J132").Value = 5300
$ws.Range("K132").Value = 4342.6155
$ws.Range("L132").Value = 15900
$ws.Range("M132").Value = -1812.6155
$ws.Range("N132").Value = -20960

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2693.889
$ws.Range("I7").Value = 2290.4
$ws.Range("J7").Value = 3198.25
$ws.Range("K7").Value = 2290.4
$ws.Range("L7").Value = 3198.25
$ws.Range("M7").Value = -2178.4
$ws.Range("N7").Value = -3422.25
# Row 82
$ws.Range("H82").Value = 6087.5
$ws.Range("I82").Value = 13300
$ws.Range("J82").Value = 3683.3333
$ws.Range("K82").Value = 13300
$ws.Range("L82").Value = 3683.3333
$ws.Range("M82").Value = -12939
$ws.Range("N82").Value = -4405.3333
# Row 85
$ws.Range("H85").Value = 6087.5
$ws.Range("I85").Value = 13300
$ws.Range("J85").Value = 3683.3333
$ws.Range("K85").Value = 13300
$ws.Range("L85").Value = 3683.3333
$ws.Range("M85").Value = -12052
$ws.Range("N85").Value = -6179.3333
# Row 100
$ws.Range("H100").Value = 1544.5
$ws.Range("I100").Value = 1544.5
$ws.Range("K100").Value = 1544.5
$ws.Range("M100").Value = -1003.5
# Row 126
$ws.Range("H126").Value = 2693.889
$ws.Range("I126").Value = 2290.4
$ws.Range("J126").Value = 3198.25
$ws.Range("K126").Value = 6871.200000000001
$ws.Range("L126").Value = 9594.75
$ws.Range("M126").Value = -4401.200000000001
$ws.Range("N126").Value = -14534.75
# Row 132
$ws.Range("H132").Value = 4789.846
$ws.Range("I132").Value = 4475.4443
$ws.Range("J132").Value = 5497.25
$ws.Range("K132").Value = 13426.3329
$ws.Range("L132").Value = 16491.75
$ws.Range("M132").Value = -10896.3329
$ws.Range("N132").Value = -21551.75
# Row 134
$ws.Range("H134").Value = 91997.92999999999
$ws.Range("J134").Value = 91997.92999999999
$ws.Range("L134").Value = 91997.92999999999
$ws.Range("N134").Value = -102137.93
# Row 136
$ws.Range("H136").Value = 2426.3333
$ws.Range("I136").Value = 2444.5
$ws.Range("K136").Value = 7333.5
$ws.Range("M136").Value = -4783.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1634.5454
$ws.Range("I122").Value = 1634.5454
$ws.Range("K122").Value = 4903.6362
$ws.Range("M122").Value = -2453.6362
# Row 124
$ws.Range("H124").Value = 23969.666
$ws.Range("J124").Value = 23969.666
$ws.Range("L124").Value = 23969.666
$ws.Range("N124").Value = -33789.666
# Row 132
$ws.Range("H132").Value = 2009.6666
$ws.Range("I132").Value = 2009.6666
$ws.Range("K132").Value = 6028.9998
$ws.Range("M132").Value = -3498.9998
